$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue ($ws.Cells.Item(2, 4)) "307.03"
Set-TextValue ($ws.Cells.Item(2, 5)) "-1.31%"
Set-TextValue ($ws.Cells.Item(2, 7)) "9"

# Row 3
Set-TextValue ($ws.Cells.Item(3, 4)) "37.11"
Set-TextValue ($ws.Cells.Item(3, 5)) "-1.55%"
Set-TextValue ($ws.Cells.Item(3, 7)) "9"

# Row 4
Set-TextValue ($ws.Cells.Item(4, 4)) "5.120"
Set-TextValue ($ws.Cells.Item(4, 5)) "0.52%"
Set-TextValue ($ws.Cells.Item(4, 7)) "9"

# Row 5
Set-TextValue ($ws.Cells.Item(5, 4)) "0.07780"
Set-TextValue ($ws.Cells.Item(5, 5)) "0.09%"
Set-TextValue ($ws.Cells.Item(5, 7)) "9"

# Row 6
Set-TextValue ($ws.Cells.Item(6, 4)) "8.240"
Set-TextValue ($ws.Cells.Item(6, 5)) "0.48%"
Set-TextValue ($ws.Cells.Item(6, 7)) "9"

# Row 7
Set-TextValue ($ws.Cells.Item(7, 4)) "1.881"
Set-TextValue ($ws.Cells.Item(7, 5)) "-1.04%"
Set-TextValue ($ws.Cells.Item(7, 7)) "9"

# Row 8
Set-TextValue ($ws.Cells.Item(8, 4)) "2.992"
Set-TextValue ($ws.Cells.Item(8, 5)) "3.63%"
Set-TextValue ($ws.Cells.Item(8, 7)) "9"

# Row 9
Set-TextValue ($ws.Cells.Item(9, 4)) "0.9277"
Set-TextValue ($ws.Cells.Item(9, 5)) "0.69%"
Set-TextValue ($ws.Cells.Item(9, 7)) "9"

# Row 10
Set-TextValue ($ws.Cells.Item(10, 4)) "0.1087"
Set-TextValue ($ws.Cells.Item(10, 5)) "-10.56%"
Set-TextValue ($ws.Cells.Item(10, 7)) "9"

# Row 11
Set-TextValue ($ws.Cells.Item(11, 4)) "0.1910"
Set-TextValue ($ws.Cells.Item(11, 5)) "-0.58%"
Set-TextValue ($ws.Cells.Item(11, 7)) "9"

# Row 12
Set-TextValue ($ws.Cells.Item(12, 4)) "0.08934"
Set-TextValue ($ws.Cells.Item(12, 5)) "-4.31%"
Set-TextValue ($ws.Cells.Item(12, 7)) "9"

# Row 13
Set-TextValue ($ws.Cells.Item(13, 4)) "0.03332"
Set-TextValue ($ws.Cells.Item(13, 5)) "-2.58%"
Set-TextValue ($ws.Cells.Item(13, 7)) "9"

# Row 14
Set-TextValue ($ws.Cells.Item(14, 4)) "0.09586"
Set-TextValue ($ws.Cells.Item(14, 5)) "-1.01%"
Set-TextValue ($ws.Cells.Item(14, 7)) "9"

# Row 15
Set-TextValue ($ws.Cells.Item(15, 4)) "0.001386"
Set-TextValue ($ws.Cells.Item(15, 5)) "1.41%"
Set-TextValue ($ws.Cells.Item(15, 7)) "9"

# Row 16
Set-TextValue ($ws.Cells.Item(16, 4)) "0.005652"
Set-TextValue ($ws.Cells.Item(16, 5)) "-5.17%"
Set-TextValue ($ws.Cells.Item(16, 7)) "9"

# Row 17
Set-TextValue ($ws.Cells.Item(17, 4)) "3.536"
Set-TextValue ($ws.Cells.Item(17, 5)) "-0.37%"
Set-TextValue ($ws.Cells.Item(17, 7)) "9"

# Row 18
Set-TextValue ($ws.Cells.Item(18, 5)) "1.56%"
Set-TextValue ($ws.Cells.Item(18, 7)) "9"

# Row 19
Set-TextValue ($ws.Cells.Item(19, 4)) "0.3368"
Set-TextValue ($ws.Cells.Item(19, 5)) "-0.93%"
Set-TextValue ($ws.Cells.Item(19, 7)) "9"

# Row 20
Set-TextValue ($ws.Cells.Item(20, 4)) "6.299"
Set-TextValue ($ws.Cells.Item(20, 5)) "19.35%"
Set-TextValue ($ws.Cells.Item(20, 7)) "9"

# Row 21
Set-TextValue ($ws.Cells.Item(21, 5)) "-1.73%"
Set-TextValue ($ws.Cells.Item(21, 7)) "9"

# Row 22
Set-TextValue ($ws.Cells.Item(22, 4)) "0.2508"
Set-TextValue ($ws.Cells.Item(22, 5)) "-3.20%"
Set-TextValue ($ws.Cells.Item(22, 7)) "9"

# Row 23
Set-TextValue ($ws.Cells.Item(23, 4)) "0.04375"
Set-TextValue ($ws.Cells.Item(23, 5)) "0.65%"
Set-TextValue ($ws.Cells.Item(23, 7)) "9"

# Row 24
Set-TextValue ($ws.Cells.Item(24, 4)) "0.001193"
Set-TextValue ($ws.Cells.Item(24, 5)) "-1.61%"
Set-TextValue ($ws.Cells.Item(24, 7)) "9"

# Row 25
Set-TextValue ($ws.Cells.Item(25, 4)) "0.004242"
Set-TextValue ($ws.Cells.Item(25, 5)) "-0.35%"
Set-TextValue ($ws.Cells.Item(25, 7)) "9"

# Row 26
Set-TextValue ($ws.Cells.Item(26, 4)) "0.0001304"
Set-TextValue ($ws.Cells.Item(26, 5)) "0.33%"
Set-TextValue ($ws.Cells.Item(26, 7)) "9"

# Row 27
Set-TextValue ($ws.Cells.Item(27, 7)) "9"

# Row 28
Set-TextValue ($ws.Cells.Item(28, 7)) "9"

# Row 29
Set-TextValue ($ws.Cells.Item(29, 7)) "9"

# Row 30
Set-TextValue ($ws.Cells.Item(30, 7)) "9"

# Row 31
Set-TextValue ($ws.Cells.Item(31, 7)) "9"

# Row 32
Set-TextValue ($ws.Cells.Item(32, 7)) "9"

# Row 33
Set-TextValue ($ws.Cells.Item(33, 7)) "9"

# Row 34
Set-TextValue ($ws.Cells.Item(34, 7)) "9"

# Row 35
Set-TextValue ($ws.Cells.Item(35, 7)) "9"

# Row 36
Set-TextValue ($ws.Cells.Item(36, 7)) "9"

# Row 37
Set-TextValue ($ws.Cells.Item(37, 7)) "9"

# Row 38
Set-TextValue ($ws.Cells.Item(38, 7)) "9"

# Row 39
Set-TextValue ($ws.Cells.Item(39, 4)) "0.02154"
Set-TextValue ($ws.Cells.Item(39, 5)) "2.38%"
Set-TextValue ($ws.Cells.Item(39, 7)) "9"

# Row 40
Set-TextValue ($ws.Cells.Item(40, 4)) "0.05017"
Set-TextValue ($ws.Cells.Item(40, 5)) "-0.76%"
Set-TextValue ($ws.Cells.Item(40, 7)) "9"

# Row 41
Set-TextValue ($ws.Cells.Item(41, 4)) "0.007472"
Set-TextValue ($ws.Cells.Item(41, 5)) "-2.60%"
Set-TextValue ($ws.Cells.Item(41, 7)) "9"

# Row 42
Set-TextValue ($ws.Cells.Item(42, 4)) "0.1346"
Set-TextValue ($ws.Cells.Item(42, 7)) "9"

# Row 43
Set-TextValue ($ws.Cells.Item(43, 4)) "0.008698"
Set-TextValue ($ws.Cells.Item(43, 5)) "-11.56%"
Set-TextValue ($ws.Cells.Item(43, 7)) "9"

# Row 44
Set-TextValue ($ws.Cells.Item(44, 5)) "2.68%"
Set-TextValue ($ws.Cells.Item(44, 7)) "9"

# Row 45
Set-TextValue ($ws.Cells.Item(45, 4)) "0.007994"
Set-TextValue ($ws.Cells.Item(45, 5)) "-9.91%"
Set-TextValue ($ws.Cells.Item(45, 7)) "9"

# Row 46
Set-TextValue ($ws.Cells.Item(46, 5)) "-1.48%"
Set-TextValue ($ws.Cells.Item(46, 7)) "9"

# Row 47
Set-TextValue ($ws.Cells.Item(47, 5)) "0.30%"
Set-TextValue ($ws.Cells.Item(47, 7)) "9"

# Row 48
Set-TextValue ($ws.Cells.Item(48, 5)) "-1.97%"
Set-TextValue ($ws.Cells.Item(48, 7)) "9"

# Row 49
Set-TextValue ($ws.Cells.Item(49, 5)) "-16.41%"
Set-TextValue ($ws.Cells.Item(49, 7)) "9"

# Row 50
Set-TextValue ($ws.Cells.Item(50, 5)) "0.30%"
Set-TextValue ($ws.Cells.Item(50, 7)) "9"

# Row 51
Set-TextValue ($ws.Cells.Item(51, 5)) "0.30%"
Set-TextValue ($ws.Cells.Item(51, 7)) "9"
